$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures (row 11 total mora, row 13 counts) ---
$ws.Range("E11").Value = 98102
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 4

# --- Copy the "last row" border formatting (currently on row 25) onto row 20 ---
# Row 20 will become the final data row after the old rows are removed below.
$ws.Range("B25:J25").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rewrite the 5 remaining employee records (new data set, part 1) ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1020808616"
$ws.Range("D16").Value = "MATEO ORTIZ AVILA"
$ws.Range("E16").Value = "1803"
$ws.Range("F16").Value = 27083
$ws.Range("G16").Value = 993642

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143368014"
$ws.Range("D17").Value = "DIOMEDES DE JESUS RENTERIA BRID"
$ws.Range("E17").Value = "1910"
$ws.Range("F17").Value = 42510
$ws.Range("G17").Value = 1356042

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047474556"
$ws.Range("D18").Value = "OSWALDO ENRIQUE OROZCO GELIS"
$ws.Range("E18").Value = "1910"
$ws.Range("F18").Value = 16562
$ws.Range("G18").Value = 828117

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73158962"
$ws.Range("D19").Value = "FERNANDO BERNAL DUQUE"
$ws.Range("E19").Value = "2309"
$ws.Range("F19").Value = 1547
$ws.Range("G19").Value = 1160000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143390596"
$ws.Range("D20").Value = "KATHIA MARGARITA MARQUEZ CUESTA"
$ws.Range("E20").Value = "2403"
$ws.Range("F20").Value = 10400
$ws.Range("G20").Value = 1300000

# --- Remove the now-obsolete employee rows (21-24) and the old row 25 whose
#     formatting we already harvested onto row 20 ---
$ws.Rows("21:25").Delete()

# --- Column D can shrink now that the longest remaining name is shorter ---
$ws.Columns("D:D").AutoFit()
